$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add dates to the "Date" column (E). Rows 2-3 use a short "d-mmm" format,
# rows 5-21 (excluding rows 4 and 11, which keep their placeholder text)
# use "mm-dd-yy".
$ws.Range("E2:E3").NumberFormat = "d-mmm"
$ws.Range("E2").Value = 43195
$ws.Range("E3").Value = 43195

$ws.Range("E5:E10,E12:E21").NumberFormat = "mm-dd-yy"
$ws.Range("E5").Value = 43197
$ws.Range("E6").Value = 43197
$ws.Range("E7").Value = 43197
$ws.Range("E8").Value = 43205
$ws.Range("E9").Value = 43205
$ws.Range("E10").Value = 43205
$ws.Range("E12").Value = 43206
$ws.Range("E13").Value = 43206
$ws.Range("E14").Value = 43205
$ws.Range("E15").Value = 43204
$ws.Range("E16").Value = 43203
$ws.Range("E17").Value = 43203
$ws.Range("E18").Value = 43204
$ws.Range("E19").Value = 43200
$ws.Range("E20").Value = 43204
$ws.Range("E21").Value = 43204

# Widen column E now that it holds full dates instead of the "xxxx" placeholder.
$ws.Columns("E").ColumnWidth = 8.3

# Update the saved view: scroll so column D is left-most and select E22
# (one row below the data, mirroring where the editor left off).
$ws.Range("E22").Select() | Out-Null
